$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds data rows 2..205 (feature index 0..203 in col A,
# value in col B). Append 12 more rows (206..217) continuing the sequence
# (feature index 204..215), reusing the same formatting (bold/centered/
# bordered style) as the existing data rows in column A.

$ws.Range("A205").Copy() | Out-Null
$ws.Range("A206:A217").PasteSpecial(-4122) | Out-Null

$aValues = @(204, 205, 206, 207, 208, 209, 210, 211, 212, 213, 214, 215)
$bValues = @(
    [double]"-3.700743415417188E-17",
    [double]"-4.037174635000569E-17",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"0",
    [double]"1.480297366166875E-16",
    [double]"0",
    [double]"0",
    [double]"1.480297366166875E-16",
    [double]"0",
    [double]"0"
)

for ($i = 0; $i -lt 12; $i++) {
    $row = 206 + $i
    $ws.Cells.Item($row, 1).Value = $aValues[$i]
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
}

Write-Host "Appended rows 206-217 to sheet1"
